$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77 is a duplicate entry of the "Polaar Techno 2025 Opening" event
# (same date/venue/city/link as the entry later in the list). Remove the
# entire row and shift everything below it up by one.
$ws.Rows.Item(77).Delete()
